# Adding the changes we made on may 9th
#
# The "ax,ay,az,gx,gy,gz" sensor readings (columns C-H) for every data row were
# regenerated:
#   - rows 4-23 now hold the sensor readings that used to live in rows 2-21
#     (i.e. the data was shifted down by two rows),
#   - rows 2-3 got brand-new sensor readings (prepended),
#   - rows 24-31 are brand-new rows appended at the end, each with its own new
#     sensor readings.
# The "timestamp" (A) and "label" (B) columns simply follow the row position
# (timestamp = 100 * (row-2), label = "struggle") and are unaffected aside
# from being filled in for the newly appended rows 22-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C,D,E,F,G,H values for every data row, indexed from row 2 to row 31.
$sensorData = @(
    @(0.04686117172241161, -0.03032520040869657, 0.02842492796480637, 0.0167987942695617, 0.1464549452066421, 0.5288565754890442),
    @(0.01706874370574948, 0.07629761099815366, -0.1732735317200423, -0.2842050492763519, -2.60824179649353, -0.3984368443489074),
    @(-0.2285218238830561, -0.162128377705812, 0.1944056674838058, -0.5183191895484924, -3.74109148979187, 0.1130100712180137),
    @(0.1096982955932607, -0.9903637319803219, 1.168615095317362, -0.485485166311264, -1.733940958976746, 0.1319468915462494),
    @(0.9691734313964835, 1.9581446647644, -7.859957635402668, -0.8013024926185608, -0.9908234477043152, -0.1090394482016563),
    @(0.1940504908561733, 0.2241474241018375, -5.310114443302172, -0.1398881375789642, 0.2066251635551452, -0.2506074607372284),
    @(-1.58094030618667, -1.941463194787501, -1.144656896591187, 1.494939923286438, 3.902817964553833, -0.2788599729537964),
    @(-1.82746779918671, -0.9869420528411925, -0.4169094562530551, 1.451415777206421, 5.435784339904785, 0.0923933684825897),
    @(-1.769100129604333, 0.7612166106700828, 0.9901726245880039, 0.4634940028190613, 1.852907061576844, 0.299934834241867),
    @(-6.232726603746395, 2.152586340904233, 3.093920886516565, 0.1380555480718612, 0.99876469373703, -0.2206750661134719),
    @(-6.777266740798977, 1.212506279349337, 2.062903136014953, -0.1919644027948379, 0.2603813111782074, -0.5108360648155212),
    @(3.723967790603678, -2.365190558135525, -3.315693378448508, -0.1438587605953216, -1.665524125099182, 0.2296853363513946),
    @(-2.335351765155792, -0.1078430823981759, -0.2542362809181207, -1.088409185409546, -4.648226261138916, -0.946077525615692),
    @(-2.176319241523743, -0.169533859938384, -1.148270010948179, -0.942259669303894, -2.754544019699097, 0.0618501044809818),
    @(-1.89803910255433, 1.068912088871008, -3.883550018072123, -0.7225008606910706, -2.54196310043335, -0.5018258094787598),
    @(-0.4242868423461992, -0.01184737682341952, -4.823844432830806, -0.4671591818332672, -1.45172131061554, -0.3419318199157715),
    @(-1.959070533514023, 0.853676103055477, 3.545509129762649, 0.116827979683876, 0.5007568001747131, 0.3149010241031647),
    @(-0.4488477408886036, 2.171596139669413, 0.5545313358307062, 1.036791086196899, 4.31255578994751, -0.2199114859104156),
    @(-1.753711834549896, 1.106546383351094, 0.5750183463096592, 1.680184841156006, 4.684267520904541, -0.0536034256219863),
    @(-6.07151210308075, 1.909203633666039, 4.961617231369019, 0.2167044430971145, 1.493718266487122, 0.555429220199585),
    @(-7.940525859594342, 0.5168425776064485, 1.82176046073439, -0.7371616363525391, 0.4766276180744171, -0.2034181207418441),
    @(-0.5248832702637101, -2.431172959506501, -2.120043188333502, -0.1747074574232101, -0.5407684445381165, -0.012980886735022),
    @(-0.4282075166702048, -0.7951091900468019, -1.767483308911331, -0.2935207486152649, -1.132696866989136, 0.0025961773935705),
    @(-2.543609619140639, 0.4462372660637008, -1.760738492012012, -0.1869247704744339, -0.6145304441452026, -0.2397646158933639),
    @(0.2467263936996389, -0.1943315342068692, -2.037849001586441, -0.1950187236070633, 0.0448985956609249, -0.3084869384765625),
    @(0.9171624183654843, 0.3103487230837345, 0.1225722581148094, 0.0444404482841491, 0.1437060534954071, -0.3220787048339844),
    @(-0.926007807254792, 0.2953229788690807, 2.136403992772098, 0.2638937830924988, 0.3306308090686798, -0.3280346393585205),
    @(-0.3645055294036889, 0.3592699170112605, 1.912454850971703, -0.0166460778564214, 0.1469130963087082, 0.2872593700885772),
    @(0.7052955031394972, 0.3671576231718077, -0.5781752001494165, -0.1206458881497383, 0.0920879393815994, 0.0424551330506801),
    @(-0.5357744693756123, -0.6449819654226362, -0.8872665241360714, 0.0326812900602817, -0.0403171069920063, 0.0131336031481623)
)

$r = 2
foreach ($rowVals in $sensorData) {
    # Fill timestamp/label for the brand-new rows appended past the original
    # 20 data rows (rows 22-31); existing rows 2-21 already have these set.
    if ($r -gt 21) {
        $ws.Cells.Item($r, 1).Value = 100 * ($r - 2)
        $ws.Cells.Item($r, 2).Value = "struggle"
    }

    $c = 3
    foreach ($v in $rowVals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c++
    }
    $r++
}
